$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4115.7915
$ws.Range("I62").Value = 4342.75
$ws.Range("J62").Value = 3661.875
$ws.Range("K62").Value = 4342.75
$ws.Range("L62").Value = 3661.875
$ws.Range("M62").Value = -3718.75
$ws.Range("N62").Value = -4909.875
$ws.Range("H65").Value = 4115.7915
$ws.Range("I65").Value = 4342.75
$ws.Range("J65").Value = 3661.875
$ws.Range("K65").Value = 21713.75
$ws.Range("L65").Value = 18309.375
$ws.Range("M65").Value = -18593.75
$ws.Range("N65").Value = -24549.375
$ws.Range("H112").Value = 3179.5
$ws.Range("J112").Value = 4499.3335
$ws.Range("L112").Value = 13498.0005
$ws.Range("N112").Value = -15714.0005
$ws.Range("H138").Value = 5309.691
$ws.Range("J138").Value = 5250.854
$ws.Range("L138").Value = 15752.562
$ws.Range("N138").Value = -26032.562

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121:N121").ClearContents()
$ws.Range("H122:N122").ClearContents()
$ws.Range("H123:N123").ClearContents()
$ws.Range("H124:N124").ClearContents()
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:N127").ClearContents()
$ws.Range("H128:N128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H135:N135").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("H119").Value = 250000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 250000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 250000
$ws.Range("N119").Value = -259676
$ws.Range("H120").Value = 190000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 190000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 190000
$ws.Range("N120").Value = -199676
$ws.Range("H122").Value = 250000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 250000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 250000
$ws.Range("N122").Value = -259800
$ws.Range("H123").Value = 59999
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 59999
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 59999
$ws.Range("N123").Value = -69799
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("H125").Value = 100000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 100000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840
$ws.Range("H126").Value = 69997
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 69997
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 69997
$ws.Range("N126").Value = -79877
$ws.Range("H127").Value = 69999
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 69999
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 69999
$ws.Range("N127").Value = -79919
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("H129").Value = 79999
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 79999
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 79999
$ws.Range("N129").Value = -89999
$ws.Range("H130").Value = 79999
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 79999
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 79999
$ws.Range("N130").Value = -90039
$ws.Range("H131").Value = 69999
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 69999
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 69999
$ws.Range("N131").Value = -80079
$ws.Range("H132").Value = 149999
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 149999
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 149999
$ws.Range("N132").Value = -160119
$ws.Range("H133").Value = 99999
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 99999
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 99999
$ws.Range("N133").Value = -110119
$ws.Range("H134").Value = 3238.9614
$ws.Range("I134").Value = 3141.3044
$ws.Range("J134").Value = 3987.6667
$ws.Range("K134").Value = 9423.913199999999
$ws.Range("L134").Value = 11963.0001
$ws.Range("M134").Value = -6888.913199999999
$ws.Range("N134").Value = -17033.0001
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H139").Value = 69748.25
$ws.Range("I139").Value = 70000
$ws.Range("J139").Value = 69664.336
$ws.Range("K139").Value = 70000
$ws.Range("L139").Value = 69664.336
$ws.Range("M139").Value = -64860
$ws.Range("N139").Value = -79944.336
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 79998.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 79998.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 79998.5
$ws.Range("N141").Value = -90358.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4487.4614
$ws.Range("I86").Value = 3786.1667
$ws.Range("K86").Value = 3786.1667
$ws.Range("M86").Value = -2663.1667
$ws.Range("H87").Value = 49499.5
$ws.Range("J87").Value = 49499.5
$ws.Range("L87").Value = 49499.5
$ws.Range("N87").Value = -51871.5
$ws.Range("H89").Value = 4487.4614
$ws.Range("I89").Value = 3786.1667
$ws.Range("K89").Value = 18930.8335
$ws.Range("M89").Value = -13314.8335
$ws.Range("H90").Value = 49499.5
$ws.Range("J90").Value = 49499.5
$ws.Range("L90").Value = 148498.5
$ws.Range("N90").Value = -160354.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120:N120").ClearContents()
$ws.Range("H121:N121").ClearContents()
$ws.Range("H122:N122").ClearContents()
$ws.Range("H123:N123").ClearContents()
$ws.Range("H124:N124").ClearContents()
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:N127").ClearContents()
$ws.Range("H128:N128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H136:N136").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:N127").ClearContents()
$ws.Range("H128:N128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H135:N135").ClearContents()
$ws.Range("H136:N136").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 35890.55
$ws.Range("I132").Value = 37065.25
$ws.Range("K132").Value = 111195.75
$ws.Range("M132").Value = -108665.75
